# Daily attendance processing - reorder "Recorded By" (column G) values so
# that "2025/2026" is listed first among the comma-separated entries.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    3  = "2025/2026, 2022/2023"
    22 = "2025/2026, 2024/2025"
    23 = "2025/2026, 2023/2024, 2022/2023"
    24 = "2025/2026, neveen.nashaat@med.asu.edu.eg"
    27 = "2025/2026, neveen.nashaat@med.asu.edu.eg"
    28 = "2025/2026, neveen.nashaat@med.asu.edu.eg"
    31 = "2025/2026, 2022/2023"
    50 = "2025/2026, 2024/2025"
    51 = "2025/2026, 2023/2024, 2022/2023"
    52 = "2025/2026, neveen.nashaat@med.asu.edu.eg"
    55 = "2025/2026, neveen.nashaat@med.asu.edu.eg"
    56 = "2025/2026, neveen.nashaat@med.asu.edu.eg"
}

foreach ($row in $updates.Keys) {
    $ws.Range("G$row").Value = $updates[$row]
}
